# Weekly price-list update: a new weekly observation (Fecha serial 44644,
# i.e. 2022-03-24) is inserted as a new data row right above the existing
# row 111, pushing all subsequent rows (old 111-131) down by one (new
# 112-132). The new row repeats the market/category/quality metadata that
# is constant across this block and carries its own Volumen/Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 111 (Excel shifts 111..131 down to 112..132).
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new weekly record.
$ws.Cells.Item(111, 1).Value = 8
$ws.Cells.Item(111, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(111, 3).Value = "Coquimbo"
$ws.Cells.Item(111, 4).Value = 44644
$ws.Cells.Item(111, 5).Value = 4
$ws.Cells.Item(111, 6).Value = 100112040
$ws.Cells.Item(111, 7).Value = "Cilantro"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 2360
$ws.Cells.Item(111, 11).Value = 2300
$ws.Cells.Item(111, 12).Value = 2500
$ws.Cells.Item(111, 13).Value = 2400
$ws.Cells.Item(111, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(111, 15).Value = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(111, 16).Value = 1600
$ws.Cells.Item(111, 17).Value = 1.5
$ws.Cells.Item(111, 18).Value = "Hortaliza"
